$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right
#    after the Heading1 title paragraph.
# ---------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Meta description", $false, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $metaPara = $searchRange.Paragraphs(1)
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# 2. Replace the final paragraph (the italic AI image-prompt paragraph)
#    with two paragraphs:
#      - a bold paragraph repeating the page title
#      - an italic paragraph containing the text that used to be the
#        "Meta description" paragraph's content (minus the "Meta
#        description:" label)
# ---------------------------------------------------------------------
$imgSearchRange = $d.Content
$imgFound = $imgSearchRange.Find.Execute("Create an image to represent the game",
                                          $false, $false, $false, $false, $false,
                                          $true, 1, $false, "", 0)
if ($imgFound) {
    $targetPara = $imgSearchRange.Paragraphs(1)
} else {
    $targetPara = $d.Paragraphs($d.Paragraphs.Count)
}

# Target the paragraph's content without its trailing paragraph mark so the
# inserted XML below replaces the paragraph's content in place rather than
# leaving a stray empty paragraph behind.
$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)

$replacementXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +
    '<w:t>Play Dragon Guard Jackpot Dash for Free - Game Review</w:t></w:r></w:p>' +
    '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr>' +
    '<w:t>Our review of Dragon Guard Jackpot Dash - a slot game with free spins, high variability, and RTP. Play for free and learn about its winning potential and design.</w:t></w:r></w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$targetRange.InsertXML($replacementXml)
